$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset formatting on A1 BEFORE writing the new value (avoids spurious font/style creation)
$ws.Range("A1").ClearFormats()

# New consolidated questions text (pretty-printed JSON-ish payload) goes into A1
$newText = @'
questions = [
    {
        "title": "You have a network diagram made in Visio. It used to fit one A4 size page, but you expanded it by adding new objects. Now it requires two pages to print. You decided to print it on one A3 size page instead.What should you do?",
        "ques_type": 2,
        "options": [
            "Change the paper size in the Print Setup section of the Page Setup window.",
            "Change the page size in the Page Size section of the Page Setup window.",
            "Change the Zoom configuration in the View ribbon.",
            "Change the Printer Properties configuration in the Print menu."
        ],
        "score": "Change the paper size in the Print Setup section of the Page Setup window."
    },
    {
        "title": "You have a network diagram, as shown below. You need to select both servers without selecting other objects.What is the easiest way to achieve that?",
        "ques_type": 2,
        "options": [
            "Navigate to Home &gt Select &gt Area Select and drag a rectangle around both servers.",
            "Go to Home &gt Select &gt Select by type and specify Networking &gt Server as a type.",
            "Select Home &gt Pointer Tool. Press Ctrl and select both servers using the mouse.",
            "Navigate to Home &gt Select &gt Lasso Select and drag a freeform lasso around the both servers."
        ],
        "score": "Select Home &gt Pointer Tool. Press Ctrl and select both servers using the mouse."
    },
    {
        "title": "You are developing a new business process for the HR team. The process is about the procedure for onboarding new employees, and it will contain multiple steps (signing of contract, assigning of mandatory training, creation of account for IT resources, provision of workplace, etc.). You want to visualize this process using Microsoft Visio diagram.Which predefined diagram template should you use?",
        "ques_type": 2,
        "options": [
            "Gantt Chart",
            "General",
            "Organization chart",
            "Flowchart"
        ],
        "score": "Flowchart"
    },
    {
        "title": "You want to visualize the organization chart of your company using Microsoft Visio. To avoid manual work, you decided to export the data from the HR tool your company uses. The tool allows for exporting data as a text file with separated values. The delimiter used for separation can be customized.Which of the following delimiters are supported by Visio?",
        "ques_type": 15,
        "options": [
            "Space",
            "Colon",
            "Semicolon",
            "Comma",
            "Tab"
        ],
        "score": [
            "Comma",
            "Tab"
        ]
    }
]
'@
$ws.Range("A1").Value = $newText

# Row 1 no longer needs the auto-sized height from the long string - restore auto height
$ws.Rows(1).AutoFit()

# Remove the now-redundant second row (its text moved into A1)
$ws.Range("A2").EntireRow.Delete()

